# Scheduled-runner market data refresh for Sargatanas_Profits workbook
# Updates currentAveragePrice / derived profit columns (H:N) per sheet/row
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H9").Value = 2446.3333
$ws.Range("I9").Value = 1899.6666
$ws.Range("J9").Value = 2993
$ws.Range("K9").Value = 1899.6666
$ws.Range("L9").Value = 2993
$ws.Range("M9").Value = -1730.6666
$ws.Range("N9").Value = -3331

$ws.Range("H17").Value = 365.2414
$ws.Range("J17").Value = 400.08
$ws.Range("L17").Value = 1200.24
$ws.Range("N17").Value = -1536.24

$ws.Range("H40").Value = 3345318.5
$ws.Range("I40").Value = 18500.334
$ws.Range("K40").Value = 18500.334
$ws.Range("M40").Value = -18325.334

$ws.Range("H41").Value = 15625154
$ws.Range("I41").Value = 20833514
$ws.Range("J41").Value = 75.5
$ws.Range("K41").Value = 20833514
$ws.Range("L41").Value = 75.5
$ws.Range("M41").Value = -20833074
$ws.Range("N41").Value = -955.5

$ws.Range("H58").Value = 9743.6
$ws.Range("I58").Value = 1500
$ws.Range("J58").Value = 11804.5
$ws.Range("K58").Value = 4500
$ws.Range("L58").Value = 35413.5
$ws.Range("M58").Value = -4350
$ws.Range("N58").Value = -35713.5

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H111").Value = 11369353
$ws.Range("I111").Value = 13894988
$ws.Range("J111").Value = 3999.5
$ws.Range("K111").Value = 41684964
$ws.Range("L111").Value = 11998.5
$ws.Range("M111").Value = -41681897
$ws.Range("N111").Value = -18132.5

$ws.Range("H125").Value = 38462410
$ws.Range("I125").Value = 55556370
$ws.Range("J125").Value = 999.5
$ws.Range("K125").Value = 500007330
$ws.Range("L125").Value = 8995.5
$ws.Range("M125").Value = -500004870
$ws.Range("N125").Value = -13915.5

$ws.Range("H132").Value = 3131
$ws.Range("I132").Value = 2779.318
$ws.Range("K132").Value = 8337.954000000002
$ws.Range("M132").Value = -5807.954000000002

$ws.Range("H135").Value = 667620.5600000001
$ws.Range("J135").Value = 3500
$ws.Range("L135").Value = 31500
$ws.Range("N135").Value = -36570

$ws.Range("H137").Value = 4229
$ws.Range("I137").Value = 4661.625
$ws.Range("J137").Value = 2498.5
$ws.Range("K137").Value = 13984.875
$ws.Range("L137").Value = 7495.5
$ws.Range("M137").Value = -11434.875
$ws.Range("N137").Value = -12595.5

$ws.Range("H138").Value = 6277.1577
$ws.Range("I138").Value = 2767.7693
$ws.Range("K138").Value = 8303.3079
$ws.Range("M138").Value = -3163.3079

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 4036399.5
$ws.Range("I32").Value = 4170782.8
$ws.Range("K32").Value = 4170782.8
$ws.Range("M32").Value = -4170495.8

$ws.Range("H45").Value = 6414.4546
$ws.Range("I45").Value = 1201
$ws.Range("J45").Value = 9393.571
$ws.Range("K45").Value = 1201
$ws.Range("L45").Value = 9393.571
$ws.Range("M45").Value = -824
$ws.Range("N45").Value = -10147.571

$ws.Range("H61").Value = 10824.632
$ws.Range("I61").Value = 5139.4165
$ws.Range("K61").Value = 5139.4165
$ws.Range("M61").Value = -4927.4165

$ws.Range("H74").Value = 50923.31
$ws.Range("I74").Value = 122799.8
$ws.Range("J74").Value = 6000.5
$ws.Range("K74").Value = 122799.8
$ws.Range("L74").Value = 6000.5
$ws.Range("M74").Value = -121925.8
$ws.Range("N74").Value = -7748.5

$ws.Range("H77").Value = 50923.31
$ws.Range("I77").Value = 122799.8
$ws.Range("J77").Value = 6000.5
$ws.Range("K77").Value = 613999
$ws.Range("L77").Value = 30002.5
$ws.Range("M77").Value = -609631
$ws.Range("N77").Value = -38738.5

$ws.Range("H122").Value = 11258.464
$ws.Range("I122").Value = 13486.4
$ws.Range("K122").Value = 40459.2
$ws.Range("M122").Value = -38009.2

$ws.Range("H136").Value = 10824.632
$ws.Range("I136").Value = 5139.4165
$ws.Range("K136").Value = 15418.2495
$ws.Range("M136").Value = -12868.2495

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H20").Value = 9805628
$ws.Range("I20").Value = 18520400
$ws.Range("J20").Value = 1509.125
$ws.Range("K20").Value = 18520400
$ws.Range("L20").Value = 1509.125
$ws.Range("M20").Value = -18520153
$ws.Range("N20").Value = -2003.125

$ws.Range("H50").Value = 45469.668
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 45469.668
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 45469.668
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -46617.668

$ws.Range("H134").Value = 5877.729
$ws.Range("I134").Value = 3054.724
$ws.Range("K134").Value = 9164.172
$ws.Range("M134").Value = -6629.172

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 8479.218000000001
$ws.Range("I31").Value = 3682.8635
$ws.Range("K31").Value = 3682.8635
$ws.Range("M31").Value = -3387.8635

$ws.Range("H34").Value = 8479.218000000001
$ws.Range("I34").Value = 3682.8635
$ws.Range("K34").Value = 3682.8635
$ws.Range("M34").Value = -3480.8635

$ws.Range("H122").Value = 3729.4443
$ws.Range("J122").Value = 4121
$ws.Range("L122").Value = 12363
$ws.Range("N122").Value = -17263

$ws.Range("H132").Value = 5276.15
$ws.Range("I132").Value = 3091.2593
$ws.Range("K132").Value = 9273.777900000001
$ws.Range("M132").Value = -6743.777900000001

$ws.Range("H134").Value = 7046.423
$ws.Range("I134").Value = 3128
$ws.Range("K134").Value = 9384
$ws.Range("M134").Value = -6849

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H33").Value = 47619308
$ws.Range("J33").Value = 517.6667
$ws.Range("L33").Value = 3106.0002
$ws.Range("N33").Value = -3672.0002

$ws.Range("H139").Value = 79688.03999999999
$ws.Range("I139").Value = 134952.94
$ws.Range("J139").Value = 4326.8184
$ws.Range("K139").Value = 404858.82
$ws.Range("L139").Value = 12980.4552
$ws.Range("M139").Value = -399718.82
$ws.Range("N139").Value = -23260.4552

$ws.Range("H140").Value = 78650.96000000001
$ws.Range("I140").Value = 92182.95
$ws.Range("K140").Value = 276548.85
$ws.Range("M140").Value = -271368.85

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H57").Value = 60333.332
$ws.Range("J57").Value = 80000
$ws.Range("L57").Value = 80000
$ws.Range("N57").Value = -81640

$ws.Range("H70").Value = 8468.6
$ws.Range("I70").Value = 7452.5454
$ws.Range("K70").Value = 7452.5454
$ws.Range("M70").Value = -7182.5454

$ws.Range("H73").Value = 8468.6
$ws.Range("I73").Value = 7452.5454
$ws.Range("K73").Value = 7452.5454
$ws.Range("M73").Value = -6516.5454

$ws.Range("H113").Value = 6098.625
$ws.Range("I113").Value = 3945.4736
$ws.Range("K113").Value = 3945.4736
$ws.Range("M113").Value = -1775.4736

$ws.Range("H122").Value = 3100.282
$ws.Range("I122").Value = 2518.2593
$ws.Range("K122").Value = 7554.777900000001
$ws.Range("M122").Value = -5104.777900000001

$ws.Range("H132").Value = 4647.125
$ws.Range("I132").Value = 1609.3158
$ws.Range("K132").Value = 4827.9474
$ws.Range("M132").Value = -2297.9474

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H7").Value = 4680.12
$ws.Range("I7").Value = 3066.6
$ws.Range("J7").Value = 7100.4
$ws.Range("K7").Value = 3066.6
$ws.Range("L7").Value = 7100.4
$ws.Range("M7").Value = -2954.6
$ws.Range("N7").Value = -7324.4

$ws.Range("H26").Value = 19000
$ws.Range("J26").Value = 19000
$ws.Range("L26").Value = 19000
$ws.Range("N26").Value = -19590

$ws.Range("H40").Value = 5945.9565
$ws.Range("I40").Value = 5377.625
$ws.Range("J40").Value = 7245
$ws.Range("K40").Value = 5377.625
$ws.Range("L40").Value = 7245
$ws.Range("M40").Value = -5241.625
$ws.Range("N40").Value = -7517

$ws.Range("H126").Value = 4680.12
$ws.Range("I126").Value = 3066.6
$ws.Range("J126").Value = 7100.4
$ws.Range("K126").Value = 9199.799999999999
$ws.Range("L126").Value = 21301.2
$ws.Range("M126").Value = -6729.799999999999
$ws.Range("N126").Value = -26241.2

$ws.Range("H132").Value = 17249120
$ws.Range("I132").Value = 35717856
$ws.Range("J132").Value = 11633.4
$ws.Range("K132").Value = 107153568
$ws.Range("L132").Value = 34900.2
$ws.Range("M132").Value = -107151038
$ws.Range("N132").Value = -39960.2

$ws.Range("H136").Value = 8120.552
$ws.Range("I136").Value = 3027.5557
$ws.Range("K136").Value = 9082.667099999999
$ws.Range("M136").Value = -6532.667099999999

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H41").Value = 21999.75
$ws.Range("J41").Value = 21999.75
$ws.Range("L41").Value = 21999.75
$ws.Range("N41").Value = -22779.75

$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("N52").ClearContents()

$ws.Range("H81").Value = 16672542
$ws.Range("J81").Value = 100030000
$ws.Range("L81").Value = 200060000
$ws.Range("N81").Value = -200062122

$ws.Range("H84").Value = 16672542
$ws.Range("J84").Value = 100030000
$ws.Range("L84").Value = 1000300000
$ws.Range("N84").Value = -1000310608

$ws.Range("H107").Value = 926.64703
$ws.Range("I107").Value = 968.125
$ws.Range("J107").Value = 889.7778
$ws.Range("K107").Value = 2904.375
$ws.Range("L107").Value = 2669.3334
$ws.Range("M107").Value = -984.375
$ws.Range("N107").Value = -6509.3334

$ws.Range("H122").Value = 20163212
$ws.Range("I122").Value = 31502544
$ws.Range("K122").Value = 94507632
$ws.Range("M122").Value = -94505182

$ws.Range("H132").Value = 38509932
$ws.Range("I132").Value = 55569170
$ws.Range("K132").Value = 166707510
$ws.Range("M132").Value = -166704980

